$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value in A2
$ws.Range("A2").Value = 123

# Row 5
$ws.Range("A5").Value = "20230309Z"
$ws.Range("B5").Value = "z"
$ws.Range("C5").Value = "z"
$ws.Range("D5").Value = 0
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "1"
$ws.Range("I5").Value = 0

# Row 6
$ws.Range("A6").Value = "20230309X"
$ws.Range("B6").Value = "x"
$ws.Range("C6").Value = "x"
$ws.Range("D6").Value = 0
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "2"
$ws.Range("I6").Value = 0

# Row 7
$ws.Range("A7").Value = "20230309RA"
$ws.Range("B7").Value = "Raju Rastogi"
$ws.Range("C7").Value = "rajuShyam"
$ws.Range("D7").Value = 123123
$ws.Range("E7").Value = "abc"
$ws.Range("F7").Value = 123
$ws.Range("G7").Value = 123
$ws.Range("H7").Value = "KALTAK"
$ws.Range("I7").Value = 0

# Selection moves to A2
$ws.Range("A2").Select() | Out-Null
